$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 342
$ws.Range("F5").Value = 1530
$ws.Range("F6").Value = 729
$ws.Range("F7").Value = 648
$ws.Range("F8").Value = 1265
$ws.Range("F9").Value = 2400
$ws.Range("F10").Value = 1298
$ws.Range("F11").Value = 263
$ws.Range("F13").Value = 1934
$ws.Range("F15").Value = 5868
$ws.Range("F16").Value = 90
$ws.Range("F17").Value = 1136
$ws.Range("F18").Value = 112
$ws.Range("F19").Value = 1331
$ws.Range("F20").Value = 1288
$ws.Range("F21").Value = 1149
$ws.Range("F23").Value = 1714
$ws.Range("F24").Value = 253
$ws.Range("F25").Value = 1079
$ws.Range("F26").Value = 612
$ws.Range("F27").Value = 101
$ws.Range("F28").Value = 162
$ws.Range("F29").Value = 3485
$ws.Range("F31").Value = 1204
$ws.Range("F33").Value = 3583
$ws.Range("F34").Value = 620
$ws.Range("F35").Value = 1143
$ws.Range("F37").Value = 118
$ws.Range("F38").Value = 946
$ws.Range("F39").Value = 331
$ws.Range("F41").Value = 37
$ws.Range("G42").Value = 68
$ws.Range("F43").Value = 89
$ws.Range("F44").Value = 841
$ws.Range("F45").Value = 1029
$ws.Range("F49").Value = 49

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 256
$ws.Range("F7").Value = 428
$ws.Range("F9").Value = 480
$ws.Range("F10").Value = 13
$ws.Range("F11").Value = 357
$ws.Range("F13").Value = 128
$ws.Range("F21").Value = 565
$ws.Range("F22").Value = 205
$ws.Range("F26").Value = 64
$ws.Range("F27").Value = 64
$ws.Range("F31").Value = 23
$ws.Range("F32").Value = 128
$ws.Range("F35").Value = 31
$ws.Range("F36").Value = 106
$ws.Range("F38").Value = 169

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 3281
$ws.Range("F5").Value = 382
$ws.Range("F7").Value = 926
$ws.Range("F8").Value = 1427
$ws.Range("F10").Value = 359
$ws.Range("F11").Value = 2674
$ws.Range("F12").Value = 224
$ws.Range("F13").Value = 407
$ws.Range("F14").Value = 1064

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 382
$ws.Range("F3").Value = 926
$ws.Range("F6").Value = 342
$ws.Range("F7").Value = 359
$ws.Range("F8").Value = 2674
$ws.Range("F9").Value = 1530
$ws.Range("F10").Value = 256
$ws.Range("F11").Value = 729
$ws.Range("F12").Value = 648
$ws.Range("F13").Value = 357
$ws.Range("F14").Value = 1265
$ws.Range("F15").Value = 2400
$ws.Range("F16").Value = 224
$ws.Range("F17").Value = 1298
$ws.Range("F19").Value = 263
$ws.Range("F21").Value = 1934
$ws.Range("F23").Value = 5868
$ws.Range("F24").Value = 407
$ws.Range("F25").Value = 1136
$ws.Range("F26").Value = 112
$ws.Range("F27").Value = 1288
$ws.Range("F28").Value = 1714
$ws.Range("F29").Value = 253
$ws.Range("F30").Value = 64
$ws.Range("F31").Value = 1079
$ws.Range("F32").Value = 612
$ws.Range("F33").Value = 162
$ws.Range("F34").Value = 3489
$ws.Range("F36").Value = 1204
$ws.Range("F37").Value = 3583
$ws.Range("F38").Value = 620
$ws.Range("F40").Value = 1143
$ws.Range("F42").Value = 118
$ws.Range("F43").Value = 946
$ws.Range("G45").Value = 68
$ws.Range("F46").Value = 841
$ws.Range("F47").Value = 1029
$ws.Range("F48").Value = 169
$ws.Range("F49").Value = 169
$ws.Range("F51").Value = 49
